$wb = $excel.ActiveWorkbook

# --- Sheet "OrderList": add row 7 ---
$wsOrderList = $wb.Worksheets.Item("OrderList")
$wsOrderList.Range("A6:C6").Copy()
$wsOrderList.Range("A7:C7").PasteSpecial(-4122)
$wsOrderList.Range("A7").Value = 6
$wsOrderList.Range("B7").Value = 4
$wsOrderList.Range("C7").Value = 126

# --- Sheet "Orders": add row 20 ---
$wsOrders = $wb.Worksheets.Item("Orders")
$wsOrders.Range("A19:I19").Copy()
$wsOrders.Range("A20:I20").PasteSpecial(-4122)
$wsOrders.Range("A20").Value = 6
$wsOrders.Range("B20").Value = "rain jacket"
$wsOrders.Range("C20").Value = 10
$wsOrders.Range("D20").Value = 21
$wsOrders.Range("E20").Value = 210
$wsOrders.Range("F20").Value = 0.4
$wsOrders.Range("G20").Value = 126
$wsOrders.Range("H20").Value = "Prawnz Store"
$wsOrders.Range("I20").Value = "placed"
